$wb = $excel.ActiveWorkbook

# Rename Sheet1 to "coefs"
$wb.Worksheets.Item("Sheet1").Name = "coefs"

# Delete Sheet2 (it is empty / no longer needed)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true
